$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Statistics")
$ws2 = $wb.Worksheets.Item("Accidents")

# --- Statistics sheet: update rows 2-36 with new simulation values ---
$ws1.Cells.Item(2,1).Value = "2024-08-30 21:16:03"
$ws1.Cells.Item(2,2).Value = 33.36845952394092
$ws1.Cells.Item(2,3).Value = 5
$ws1.Cells.Item(3,1).Value = "2024-08-30 21:16:05"
$ws1.Cells.Item(3,2).Value = 36.18499795343281
$ws1.Cells.Item(3,3).Value = 8
$ws1.Cells.Item(4,1).Value = "2024-08-30 21:16:07"
$ws1.Cells.Item(4,2).Value = 38.13729141410148
$ws1.Cells.Item(4,3).Value = 10
$ws1.Cells.Item(5,1).Value = "2024-08-30 21:16:09"
$ws1.Cells.Item(5,2).Value = 38.54517290604009
$ws1.Cells.Item(5,3).Value = 13
$ws1.Cells.Item(6,1).Value = "2024-08-30 21:16:11"
$ws1.Cells.Item(6,2).Value = 40.10152733482884
$ws1.Cells.Item(6,3).Value = 17
$ws1.Cells.Item(7,1).Value = "2024-08-30 21:16:14"
$ws1.Cells.Item(7,2).Value = 37.83928087094589
$ws1.Cells.Item(7,3).Value = 19
$ws1.Cells.Item(8,1).Value = "2024-08-30 21:16:16"
$ws1.Cells.Item(8,2).Value = 33.64837695736625
$ws1.Cells.Item(8,3).Value = 21
$ws1.Cells.Item(9,1).Value = "2024-08-30 21:16:18"
$ws1.Cells.Item(9,2).Value = 33.91395395140444
$ws1.Cells.Item(9,3).Value = 26
$ws1.Cells.Item(10,1).Value = "2024-08-30 21:16:20"
$ws1.Cells.Item(10,2).Value = 28.2301595680891
$ws1.Cells.Item(10,3).Value = 26
$ws1.Cells.Item(11,1).Value = "2024-08-30 21:16:22"
$ws1.Cells.Item(11,2).Value = 29.85904993760463
$ws1.Cells.Item(11,3).Value = 28
$ws1.Cells.Item(12,1).Value = "2024-08-30 21:16:24"
$ws1.Cells.Item(12,2).Value = 29.5493732380974
$ws1.Cells.Item(12,3).Value = 29
$ws1.Cells.Item(13,1).Value = "2024-08-30 21:16:26"
$ws1.Cells.Item(13,2).Value = 27.00575810892638
$ws1.Cells.Item(13,3).Value = 33
$ws1.Cells.Item(14,1).Value = "2024-08-30 21:16:28"
$ws1.Cells.Item(14,2).Value = 25.19415498944947
$ws1.Cells.Item(14,3).Value = 35
$ws1.Cells.Item(15,1).Value = "2024-08-30 21:16:30"
$ws1.Cells.Item(15,2).Value = 23.83743746686357
$ws1.Cells.Item(15,3).Value = 37
$ws1.Cells.Item(16,1).Value = "2024-08-30 21:16:32"
$ws1.Cells.Item(16,2).Value = 22.58287766971266
$ws1.Cells.Item(16,3).Value = 38
$ws1.Cells.Item(17,1).Value = "2024-08-30 21:16:34"
$ws1.Cells.Item(17,2).Value = 21.42318548347823
$ws1.Cells.Item(17,3).Value = 37
$ws1.Cells.Item(18,1).Value = "2024-08-30 21:17:12"
$ws1.Cells.Item(18,2).Value = 19.97962162743501
$ws1.Cells.Item(18,3).Value = 38
$ws1.Cells.Item(19,1).Value = "2024-08-30 21:17:14"
$ws1.Cells.Item(19,2).Value = 16.73296032689855
$ws1.Cells.Item(19,3).Value = 38
$ws1.Cells.Item(20,1).Value = "2024-08-30 21:17:16"
$ws1.Cells.Item(20,2).Value = 17.88842700428664
$ws1.Cells.Item(20,3).Value = 41
$ws1.Cells.Item(21,1).Value = "2024-08-30 21:17:18"
$ws1.Cells.Item(21,2).Value = 19.44136268782263
$ws1.Cells.Item(21,3).Value = 40
$ws1.Cells.Item(22,1).Value = "2024-08-30 21:17:20"
$ws1.Cells.Item(22,2).Value = 17.24289317434194
$ws1.Cells.Item(22,3).Value = 39
$ws1.Cells.Item(23,1).Value = "2024-08-30 21:17:22"
$ws1.Cells.Item(23,2).Value = 17.4161599985886
$ws1.Cells.Item(23,3).Value = 40
$ws1.Cells.Item(24,1).Value = "2024-08-30 21:17:24"
$ws1.Cells.Item(24,2).Value = 15.91116227216883
$ws1.Cells.Item(24,3).Value = 39
$ws1.Cells.Item(25,1).Value = "2024-08-30 21:17:26"
$ws1.Cells.Item(25,2).Value = 16.3350609746136
$ws1.Cells.Item(25,3).Value = 38
$ws1.Cells.Item(26,1).Value = "2024-08-30 21:17:28"
$ws1.Cells.Item(26,2).Value = 17.51482579951818
$ws1.Cells.Item(26,3).Value = 39
$ws1.Cells.Item(27,1).Value = "2024-08-30 21:17:30"
$ws1.Cells.Item(27,2).Value = 17.20392362997698
$ws1.Cells.Item(27,3).Value = 38
$ws1.Cells.Item(28,1).Value = "2024-08-30 21:17:32"
$ws1.Cells.Item(28,2).Value = 17.81991207315941
$ws1.Cells.Item(28,3).Value = 39
$ws1.Cells.Item(29,1).Value = "2024-08-30 21:17:34"
$ws1.Cells.Item(29,2).Value = 15.91981948055395
$ws1.Cells.Item(29,3).Value = 38
$ws1.Cells.Item(30,1).Value = "2024-08-30 21:17:36"
$ws1.Cells.Item(30,2).Value = 14.49940314875241
$ws1.Cells.Item(30,3).Value = 38
$ws1.Cells.Item(31,1).Value = "2024-08-30 21:17:38"
$ws1.Cells.Item(31,2).Value = 15.5213523136577
$ws1.Cells.Item(31,3).Value = 39
$ws1.Cells.Item(32,1).Value = "2024-08-30 21:17:40"
$ws1.Cells.Item(32,2).Value = 17.98013377150679
$ws1.Cells.Item(32,3).Value = 38
$ws1.Cells.Item(33,1).Value = "2024-08-30 21:17:42"
$ws1.Cells.Item(33,2).Value = 17.70205052664934
$ws1.Cells.Item(33,3).Value = 39
$ws1.Cells.Item(34,1).Value = "2024-08-30 21:17:44"
$ws1.Cells.Item(34,2).Value = 17.85047680230193
$ws1.Cells.Item(34,3).Value = 39
$ws1.Cells.Item(35,1).Value = "2024-08-30 21:17:46"
$ws1.Cells.Item(35,2).Value = 15.19020394449197
$ws1.Cells.Item(35,3).Value = 39
$ws1.Cells.Item(36,1).Value = "2024-08-30 21:17:48"
$ws1.Cells.Item(36,2).Value = 15.25731561172521
$ws1.Cells.Item(36,3).Value = 40

# Remove now-unused trailing rows 37-42 (simulation window shortened)
$ws1.Rows("37:42").Delete()

# --- Accidents sheet: correct row 2 and append a new accident row 3 ---
$ws2.Cells.Item(2,1).Value = "2024-08-30 21:16:32"
$ws2.Cells.Item(2,2).Value = "Car and Truck"
$ws2.Cells.Item(2,3).Value = "28.77 and 20.30"
$ws2.Cells.Item(2,4).Value = 1

$ws2.Cells.Item(3,1).Value = "2024-08-30 21:16:34"
$ws2.Cells.Item(3,2).Value = "Car and Truck"
$ws2.Cells.Item(3,3).Value = "29.10 and 0.00"
$ws2.Cells.Item(3,4).Value = 1
